$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correcting order of columns in Table 4: swap "Dust aerosol" (B) and
# "Sea salt aerosol" (C) columns so Sea salt aerosol comes first (col B).
for ($r = 1; $r -le 8; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}
